$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add the two new task rows (rows 9 and 10), appended to sharedStrings
#        at the end, same order as the target (index 8 = "Colores en forms",
#        index 9 = "Que si cancelo imprimir..."). Column A keeps the current
#        (pre-sort) top-to-bottom order from the original list.
$ws.Range("A9").Value = "Colores en forms"
$ws.Range("A10").Value = "Que si cancelo imprimir al crear partida, me cancele todas las impresiones"

# --- 2. Fill in column B ("group/priority" order key) for every row, in the
#        CURRENT (pre-sort) row order. Row 1 is the already-done task (kept
#        pinned at the very top), rows 2-10 will be reordered by the sort
#        step below, and a stable ascending sort on these keys reproduces
#        the exact final row order from the target workbook.
$ws.Range("B1").Value = 1    # Que no me deje asociar la partida mas de una vez   -> done
$ws.Range("B2").Value = 1    # que no me deje confirmar la adquisicion...          -> done
$ws.Range("B3").Value = 2    # Documentos de rendicion, rutas, idioma e imprimir
$ws.Range("B4").Value = 4    # Poder modificar las rutas desde el sistema
$ws.Range("B5").Value = 4    # Poder modificar el mail desde el sistema
$ws.Range("B6").Value = 5    # Que todos lso form tengan singleton...
$ws.Range("B7").Value = 3    # Traduccion en los messageBox
$ws.Range("B8").Value = 3    # Traduccion en las grillas
$ws.Range("B9").Value = 5    # Colores en forms
$ws.Range("B10").Value = 1   # Que si cancelo imprimir al crear partida...        -> done

# --- 3. Highlight the two completed ("Realizado") tasks in green.
$ws.Range("A1:A2").Interior.Color = 5296274   # BGR for RGB(146, 208, 80) / FF92D050

# --- 4. Column A width, sized to fit the longest task description
#        (closest the engine's column-width rounding can get to 75.140625).
$ws.Columns("A").ColumnWidth = 74.3

# --- 5. Sort rows 2-10 by column B ascending (row 1 / the header-like first
#        task stays pinned in place), leaving a recorded sortState behind.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B10")) | Out-Null
$ws.Sort.SetRange($ws.Range("A2:B10"))
$ws.Sort.Apply()

# --- 6. View state: zoomed in, with the selection parked below the data.
$excel.ActiveWindow.Zoom = 190
$ws.Range("B12").Select()
